$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns: team record (Wins / Losses / Ties) appended after the
# existing AC (Unnamed: 28) column, i.e. AD, AE, AF.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the bold / centered / bordered header formatting used by A1:AC1.
$headerRange = $ws.Range("AD1:AF1")
$headerRange.Font.Bold = $true
$headerRange.HorizontalAlignment = -4108  # xlCenter
$headerRange.VerticalAlignment = -4160    # xlTop
$headerRange.Borders.LineStyle = 1
$headerRange.Borders.Weight = 2

# Same team record repeated for every player row (2-49).
$lastRow = 49
$ws.Range("AD2:AD$lastRow").Value = 68
$ws.Range("AE2:AE$lastRow").Value = 94
$ws.Range("AF2:AF$lastRow").Value = 0
